$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (column D) / Volume(1h) (column E) updates for this refresh, keyed by
# cell address. Values are written through a scratch cell that is explicitly
# formatted as Text and then copied in via PasteSpecial(xlPasteValues) so that
# numeric-looking strings (e.g. "1.00", "4.75") are preserved verbatim as text
# instead of being auto-coerced into numbers (which would drop trailing zeros
# or introduce floating point noise) - matching the original inline-string
# cells, and the scratch cell/style is cleaned up afterwards so no stray
# formatting or used-range growth is left behind.
$updates = @(
    @{ Cell = "D2"; Value = "56.599.79" },
    @{ Cell = "E2"; Value = "  -3.84%  " },
    @{ Cell = "D3"; Value = "2.376.69" },
    @{ Cell = "E3"; Value = "  -4.78%  " },
    @{ Cell = "E4"; Value = "  -0.17%  " },
    @{ Cell = "D5"; Value = "512.15" },
    @{ Cell = "E5"; Value = "  -4.37%  " },
    @{ Cell = "D6"; Value = "130.79" },
    @{ Cell = "E6"; Value = "  -2.84%  " },
    @{ Cell = "D7"; Value = "0.996" },
    @{ Cell = "E7"; Value = "  -0.32%  " },
    @{ Cell = "E8"; Value = "  -2.35%  " },
    @{ Cell = "D9"; Value = "2.398.61" },
    @{ Cell = "E9"; Value = "  -4.83%  " },
    @{ Cell = "D10"; Value = "0.0964" },
    @{ Cell = "E10"; Value = "  -3.21%  " },
    @{ Cell = "E11"; Value = "  -1.47%  " },
    @{ Cell = "D12"; Value = "4.75" },
    @{ Cell = "E12"; Value = "  -9.79%  " },
    @{ Cell = "D13"; Value = "0.321" },
    @{ Cell = "E13"; Value = "  -3.77%  " },
    @{ Cell = "D14"; Value = "2.799.52" },
    @{ Cell = "E14"; Value = "  -5.07%  " },
    @{ Cell = "D15"; Value = "56.532.24" },
    @{ Cell = "E15"; Value = "  -3.80%  " },
    @{ Cell = "D16"; Value = "21.66" },
    @{ Cell = "E16"; Value = "  -3.78%  " },
    @{ Cell = "E17"; Value = "  -3.41%  " },
    @{ Cell = "D18"; Value = "2.382.54" },
    @{ Cell = "E18"; Value = "  -5.50%  " },
    @{ Cell = "D19"; Value = "10.35" },
    @{ Cell = "E19"; Value = "  -3.18%  " },
    @{ Cell = "D20"; Value = "313.16" },
    @{ Cell = "E20"; Value = "  -2.66%  " },
    @{ Cell = "D21"; Value = "4.07" },
    @{ Cell = "E21"; Value = "  -4.46%  " },
    @{ Cell = "D22"; Value = "6.23" },
    @{ Cell = "E22"; Value = "  -0.51%  " },
    @{ Cell = "D23"; Value = "1.00" },
    @{ Cell = "E23"; Value = "  +0.17%  " },
    @{ Cell = "D24"; Value = "65.14" },
    @{ Cell = "E24"; Value = "  -1.01%  " },
    @{ Cell = "D25"; Value = "0.997" },
    @{ Cell = "E25"; Value = "  -0.07%  " },
    @{ Cell = "D26"; Value = "0.391" },
    @{ Cell = "E26"; Value = "  -5.09%  " },
    @{ Cell = "D27"; Value = "2.474.15" },
    @{ Cell = "E27"; Value = "  -5.72%  " },
    @{ Cell = "D28"; Value = "0.153" },
    @{ Cell = "E28"; Value = "  -4.97%  " },
    @{ Cell = "D29"; Value = "7.26" },
    @{ Cell = "E29"; Value = "  -3.49%  " },
    @{ Cell = "D30"; Value = "174.78" },
    @{ Cell = "E30"; Value = "  +1.48%  " },
    @{ Cell = "D31"; Value = "1.69" },
    @{ Cell = "E31"; Value = "  -3.07%  " },
    @{ Cell = "D32"; Value = "0.0₃0718" },
    @{ Cell = "E32"; Value = "  -6.03%  " },
    @{ Cell = "E33"; Value = "  -3.00%  " },
    @{ Cell = "D34"; Value = "1.12" },
    @{ Cell = "E34"; Value = "  -6.97%  " },
    @{ Cell = "E35"; Value = "  -0.16%  " },
    @{ Cell = "D36"; Value = "0.994" },
    @{ Cell = "E36"; Value = "  +0.06%  " },
    @{ Cell = "D37"; Value = "17.77" },
    @{ Cell = "E37"; Value = "  -2.44%  " },
    @{ Cell = "E38"; Value = "  -4.02%  " },
    @{ Cell = "E39"; Value = "  -7.12%  " },
    @{ Cell = "D40"; Value = "35.76" },
    @{ Cell = "E40"; Value = "  -2.60%  " },
    @{ Cell = "E41"; Value = "  -5.11%  " },
    @{ Cell = "D42"; Value = "0.797" },
    @{ Cell = "E42"; Value = "  -2.80%  " },
    @{ Cell = "D43"; Value = "129.54" },
    @{ Cell = "E43"; Value = "  -2.09%  " },
    @{ Cell = "D44"; Value = "3.36" },
    @{ Cell = "E44"; Value = "  -3.90%  " },
    @{ Cell = "D45"; Value = "4.92" },
    @{ Cell = "E45"; Value = "  -5.42%  " },
    @{ Cell = "D46"; Value = "256.53" },
    @{ Cell = "E46"; Value = "  -7.78%  " },
    @{ Cell = "D47"; Value = "0.575" },
    @{ Cell = "E47"; Value = "  -3.43%  " },
    @{ Cell = "D48"; Value = "0.0901" },
    @{ Cell = "E48"; Value = "  -3.70%  " },
    @{ Cell = "D49"; Value = "0.0490" },
    @{ Cell = "E49"; Value = "  -4.45%  " },
    @{ Cell = "D50"; Value = "0.0208" },
    @{ Cell = "E50"; Value = "  -5.55%  " },
    @{ Cell = "D51"; Value = "16.75" },
    @{ Cell = "E51"; Value = "  -5.77%  " }
)

$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

foreach ($item in $updates) {
    $helper.Value = $item.Value
    $helper.Copy()
    $ws.Range($item.Cell).PasteSpecial(-4163)
}

$helper.Clear()
$excel.CutCopyMode = $false
